$d = $word.ActiveDocument

# --- Paragraph: "Basic information about requirements ..." ---
# Merge the split runs (around the "be run" gramStart/gramEnd proofErr markers)
# back into a single run with identical text.
$d.Content.Find.Execute(
    "destined to be run on the laptops",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "destined to be run on the laptops", 2)

# --- Paragraph: "Resources are fairly limited, ..." ---
# Merge the split runs (around the "fairly limited" gramStart/gramEnd proofErr
# markers) back into a single run with identical text.
$d.Content.Find.Execute(
    "Resources are fairly limited, as are potential vendor options",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Resources are fairly limited, as are potential vendor options", 2)

# --- Paragraph: "Given these assumptions, ..." ---
# Replace "many laptop remain to be purchased." with "mdfgfdgfdgd" -- this
# collapses the trailing text right up to (and including) the position of the
# _GoBack bookmark, and also removes the "be purchased" gramStart/gramEnd
# proofErr markers by merging the runs.
$d.Content.Find.Execute(
    "many laptop remain to be purchased.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "mdfgfdgfdgd", 2)

# Now insert the remaining new text "any laptop remain to be purchased."
# immediately after the _GoBack bookmark, so the bookmark ends up sitting
# between the "...how mdfgfdgfdgd" run and the new "any laptop..." run,
# exactly as in the target document.
$bm = $d.Bookmarks("_GoBack")
$pos = $bm.Start
$bm.Delete()
$r = $d.Range($pos, $pos)
$r.InsertAfter("any laptop remain to be purchased.")
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos))
